# Deploy the implementation guide.
# - Update the "Date" metadata value
# - Update the "Contact" metadata value
# - Insert a new "Jurisdiction" metadata row (empty value) right after
#   "Contact", pushing "Description" and everything below it down by one
#   row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Date (row 8, column B)
$ws.Range("B8").Value = "2024-10-02T15:04:17+00:00"

# Update Contact (row 10, column B)
$ws.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# Insert a blank row after Contact (row 10) for the new Jurisdiction
# entry. This pushes Description and everything below it down by one
# row.
$ws.Rows.Item(11).Insert() | Out-Null

# The freshly inserted row gets Excel's blank default formatting, so
# copy the (now-shifted) body-row formatting from row 12 onto row 11
# to match the rest of the table (bordered, wrap-text body style).
$ws.Range("A12:B12").Copy() | Out-Null
$ws.Range("A11:B11").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
